$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.258.09'
$ws.Range("E2").Value = '  -3.05%  '

$ws.Range("D3").Value = '3.068.73'
$ws.Range("E3").Value = '  -3.92%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.82'
$ws.Range("E5").Value = '  -2.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.24'
$ws.Range("E6").Value = '  -5.74%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.066.34'
$ws.Range("E8").Value = '  -3.91%  '

$ws.Range("E9").Value = '  +3.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.22'
$ws.Range("E10").Value = '  -0.18%  '

$ws.Range("E11").Value = '  -5.46%  '

$ws.Range("E12").Value = '  -0.10%  '

$ws.Range("E13").Value = '  +1.80%  '

$ws.Range("D14").Value = '3.596.00'
$ws.Range("E14").Value = '  -3.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.09'
$ws.Range("E15").Value = '  -4.49%  '

$ws.Range("E16").Value = '  -5.67%  '

$ws.Range("D17").Value = '57.273.63'
$ws.Range("E17").Value = '  -2.99%  '

$ws.Range("D18").Value = '3.067.42'
$ws.Range("E18").Value = '  -3.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.86'
$ws.Range("E19").Value = '  -6.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.41'
$ws.Range("E20").Value = '  -5.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.82'
$ws.Range("E21").Value = '  -4.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '348.10'
$ws.Range("E22").Value = '  -4.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.13'
$ws.Range("E24").Value = '  -0.47%  '

$ws.Range("E25").Value = '  -4.61%  '

$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.165'
$ws.Range("E27").Value = '  -3.96%  '

$ws.Range("D28").Value = '0.0₃0842'
$ws.Range("E28").Value = '  -13.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.06'
$ws.Range("E30").Value = '  -7.26%  '

$ws.Range("E31").Value = '  -3.82%  '

$ws.Range("E32").Value = '  -11.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.00'
$ws.Range("E33").Value = '  -2.95%  '

$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '158.57'
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.80'
$ws.Range("E35").Value = '  -2.59%  '

$ws.Range("E36").Value = '  -8.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.96'
$ws.Range("E37").Value = '  -6.42%  '

$ws.Range("E38").Value = '  -4.02%  '

$ws.Range("E39").Value = '  -7.78%  '

$ws.Range("E40").Value = '  -3.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.57'
$ws.Range("E41").Value = '  -6.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.00'
$ws.Range("E42").Value = '  -1.96%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.687'
$ws.Range("E43").Value = '  -3.88%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.406.70'
$ws.Range("E44").Value = '  +1.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.73'
$ws.Range("E45").Value = '  -1.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("D47").Value = '3.106.47'
$ws.Range("E47").Value = '  -3.73%  '

$ws.Range("E48").Value = '  -5.29%  '

$ws.Range("E49").Value = '  -2.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.939'
$ws.Range("E50").Value = '  -9.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.28'
$ws.Range("E51").Value = '  -8.35%  '
